$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("102_1")

# Row 25: "Measures reported, House joint resolutions" - change Senate (B) value from 11 to 1
$ws.Range("B25").Value = 1

# Row 40: "Bills vetoed" - add Senate (B) value of 1
$ws.Range("B40").Value = 1
